# Apply line breaks inside four paragraphs of the "LOB1038" course sheet:
#  1) The Portuguese "Programa" paragraph (numbered items 1-9)
#  2) The English "Programa" paragraph (italic, numbered items 1-9)
#  3) The "Norma de recuperação" sentence ("O " / "(NF+RC)/2 ...")
#  4) The "Bibliografia" paragraph (one reference per line)
#
# Each insertion point is turned into a <w:br/> by searching the literal
# text spanning the join point and replacing it with the same text plus a
# "^l" (manual line break) wildcard-free Find/Replace token.

$d = $word.ActiveDocument

function Break-At($searchText, $replaceText) {
    $found = $d.Content.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)
    if (-not $found) {
        throw "Find/Replace failed for: $searchText"
    }
}

# --- Paragraph: Portuguese "Programa" -------------------------------------
Break-At "significativos.2) Construção" "significativos.^l2) Construção"
Break-At "Linearização.3) Regressão" "Linearização.^l3) Regressão"
Break-At "quadrados. 4) Cinemática" "quadrados. ^l4) Cinemática"
Break-At "Queda Livre.5) Estática" "Queda Livre.^l5) Estática"
Break-At "Material. 6) Atrito." "Material. ^l6) Atrito."
Break-At "Atrito.7) Lei de Hooke" "Atrito.^l7) Lei de Hooke"
Break-At "Young. 8) Conservação" "Young. ^l8) Conservação"
Break-At "Massa-mola.9) Choques" "Massa-mola.^l9) Choques"

# --- Paragraph: English "Programa" (italic) --------------------------------
Break-At "figures.2) Construction" "figures.^l2) Construction"
Break-At "Linearization.3) Introduction" "Linearization.^l3) Introduction"
Break-At "minimum.4) Kinematics" "minimum.^l4) Kinematics"
Break-At "Free fall.5) Statics" "Free fall.^l5) Statics"
Break-At "material point. 6) Friction." "material point. ^l6) Friction."
Break-At "Friction.7)" "Friction.^l7)"
Break-At "Young´s Modulus.8) Energy" "Young´s Modulus.^l8) Energy"
Break-At "Mass-spring system.9) Shocks." "Mass-spring system.^l9) Shocks."

# --- "Norma de recuperação" sentence ---------------------------------------
Break-At "O (NF+RC)/2" "O ^l(NF+RC)/2"

# --- Paragraph: "Bibliografia" ---------------------------------------------
Break-At "IFSC/USP.CRUZ," "IFSC/USP.^lCRUZ,"
Break-At "Curso deLaboratório:" "Curso de^lLaboratório:"
Break-At "UNICAMP (2005).NUSSENZVEIG," "UNICAMP (2005).^lNUSSENZVEIG,"
Break-At "Blucher (2008).RESNICK," "Blucher (2008).^lRESNICK,"
Break-At "LTC (2008).TIPLER," "LTC (2008).^lTIPLER,"
Break-At "LTC (2008).SEARS," "LTC (2008).^lSEARS,"
Break-At "Wesley (2009).JEWETT" "Wesley (2009).^lJEWETT"
